# Swap the two team-member names in the "Individual Contributions" table
# on the second slide (sldId 258): row 2 ("Akhil Patlori") and row 4
# ("Vikram Boppana") exchange their names.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the table shape on the slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
        break
    }
}

$tbl = $tableShape.Table

# Row 2, column 1 currently holds "Akhil Patlori" -> becomes "Vikram Boppana"
$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Vikram Boppana"

# Row 4, column 1 currently holds "Vikram Boppana" -> becomes "Akhil Patlori"
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Akhil Patlori"
